# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a handful of rows in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: sd / Statement-non-opinion -> sv / Statement-opinion
$ws.Cells.Item(13, 9).Value = "sv"
$ws.Cells.Item(13, 10).Value = "Statement-opinion"

# Row 21: aa / Agree/Accept -> sd / Statement-non-opinion
$ws.Cells.Item(21, 9).Value = "sd"
$ws.Cells.Item(21, 10).Value = "Statement-non-opinion"

# Row 22: sd / Statement-non-opinion -> sv / Statement-opinion
$ws.Cells.Item(22, 9).Value = "sv"
$ws.Cells.Item(22, 10).Value = "Statement-opinion"

# Row 28: qy / Yes-No-Question -> ba / Appreciation
$ws.Cells.Item(28, 9).Value = "ba"
$ws.Cells.Item(28, 10).Value = "Appreciation"

# Row 40: sd / Statement-non-opinion -> sv / Statement-opinion
$ws.Cells.Item(40, 9).Value = "sv"
$ws.Cells.Item(40, 10).Value = "Statement-opinion"

# Row 48: sv / Statement-opinion -> sd / Statement-non-opinion
$ws.Cells.Item(48, 9).Value = "sd"
$ws.Cells.Item(48, 10).Value = "Statement-non-opinion"
